$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "BHARANI KUMAR M"
$ws.Range("C8").Value = "/static/images/profile_photos/002/VEC-002-05-2.webp"
$ws.Range("J8").Value = "VEC-002-05-2"
$ws.Range("B8").Value = "Lab Instructor"
